# Purchasing / Operations data correction on Sheet1.
# Rows 2-6 (Round 0) and 7-11 (Round 1) previously had placeholder "?"
# values for Qlty (E), TransP_mode (J) and Trade_unit (K); also a few
# Name (B) / Country (D) typos for the Pure/Vital suppliers. This fills
# in the correct values (matching the pattern already used by the
# Round 2 / Round 3 rows further down the sheet), and corrects the
# "map" country/name data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Round 0 (rows 2-6) ---
$ws.Range("E2").Value = "High"
$ws.Range("J2").Value = "Truck"
$ws.Range("K2").Value = "Pallet"

$ws.Range("B3").Value = "WWP"
$ws.Range("E3").Value = "High"
$ws.Range("J3").Value = "Truck"
$ws.Range("K3").Value = "Pallet"

$ws.Range("J4").Value = "Boat"
$ws.Range("K4").Value = "Tank"

$ws.Range("B5").Value = "Pure"
$ws.Range("D5").Value = "Congo"
$ws.Range("E5").Value = "Middle"
$ws.Range("J5").Value = "Boat"
$ws.Range("K5").Value = "IBC"

$ws.Range("B6").Value = "Vital"
$ws.Range("D6").Value = "China"
$ws.Range("E6").Value = "Middle"
$ws.Range("J6").Value = "Boat"
$ws.Range("K6").Value = "Drum"

# --- Round 1 (rows 7-11) ---
$ws.Range("E7").Value = "High"
$ws.Range("J7").Value = "Truck"
$ws.Range("K7").Value = "Pallet"

$ws.Range("B8").Value = "WWP"
$ws.Range("E8").Value = "High"
$ws.Range("J8").Value = "Truck"
$ws.Range("K8").Value = "Pallet"

$ws.Range("J9").Value = "Boat"
$ws.Range("K9").Value = "Tank"

$ws.Range("B10").Value = "Pure"
$ws.Range("D10").Value = "Congo"
$ws.Range("E10").Value = "Middle"
$ws.Range("J10").Value = "Boat"
$ws.Range("K10").Value = "IBC"

$ws.Range("B11").Value = "Vital"
$ws.Range("D11").Value = "China"
$ws.Range("E11").Value = "High"
$ws.Range("J11").Value = "Boat"
$ws.Range("K11").Value = "Drum"

# Move the visible selection to K13, matching where the editor left off.
$ws.Range("K13").Select() | Out-Null
